$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.302.57"
$ws.Range("E2").Value = "  +2.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.86"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.98"
$ws.Range("E5").Value = "  +2.22%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.26"
$ws.Range("E9").Value = "  +2.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08006"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.008"
$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.77"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.289"
$ws.Range("E15").Value = "  +2.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.21"
$ws.Range("E16").Value = "  +2.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.300.47"
$ws.Range("E22").Value = "  +2.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.457"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.101.40"
$ws.Range("E26").Value = "  +1.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.96"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.146"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.21"
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9784"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09503"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.594"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.379"
$ws.Range("E35").Value = "  +2.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.359"
$ws.Range("E36").Value = "  +1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02276"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06098"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.436"
$ws.Range("E39").Value = "  +1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5981"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1883"
$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.41"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.289"
$ws.Range("E45").Value = "  +3.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5612"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.13"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.970"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06912"
$ws.Range("E49").Value = "  +3.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.27"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.037"
$ws.Range("E51").Value = "  +13.89%  "

Write-Output "Applied crypto price updates"